$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '244.88'
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '23.11'
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.433'
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.05969'
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.8096'
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.9225'
$cell.Style = "Normal"

$ws.Range("B9").Value = 'WazirX'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1428'
$cell.Style = "Normal"

$ws.Range("E9").Value = '8WazirXWRX'

$ws.Range("B10").Value = 'MandalaExchangeToken'

$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07428'
$cell.Style = "Normal"

$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.03394'
$cell.Style = "Normal"

$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B12").Value = 'BitrueCoin'

$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.03037'
$cell.Style = "Normal"

$ws.Range("E12").Value = '11BitrueCoinBTR'

$ws.Range("B13").Value = 'BitMartToken'

$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.09348'
$cell.Style = "Normal"

$ws.Range("E13").Value = '12BitMartTokenBMX'

$ws.Range("B14").Value = 'MCDex'

$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '3.941'
$cell.Style = "Normal"

$ws.Range("E14").Value = '13MCDexMCB'

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.001595'
$cell.Style = "Normal"

$ws.Range("E15").Value = '14BitForexTokenBF'

$ws.Range("B16").Value = 'CoinExToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.04819'
$cell.Style = "Normal"

$ws.Range("E16").Value = '15CoinExTokenCET'

$ws.Range("B17").Value = 'One'

$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.0005943'
$cell.Style = "Normal"

$ws.Range("E17").Value = '16OneONE'

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.005468'
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.004154'
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0009824'
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.00007704'
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.452'
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.186'
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.0002447'
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.03940'
$cell.Style = "Normal"

$ws.Range("B41").Value = 'BKEXToken'

$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.1073'
$cell.Style = "Normal"

$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'CEJI'

$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.002661'
$cell.Style = "Normal"

$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'KickToken'

$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.003021'
$cell.Style = "Normal"

$ws.Range("E43").Value = '42KickTokenKICK'

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.007220'
$cell.Style = "Normal"

$ws.Range("E44").Value = '43LocalTradersLCTBestin24h'

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.00005134'
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000750'
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.8554'
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.00002101'
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0002001'
$cell.Style = "Normal"
